$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.095.84'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.24%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.835.05'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.29%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9994'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.36'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.73%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6357'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.62%  '

$ws.Range("E7").Value = '  +0.17%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07558'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.30%  '

$ws.Range("E9").Value = '  +1.34%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.35'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.48%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07734'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.24%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.840.45'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.62%  '

$ws.Range("E13").Value = '  +1.15%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6737'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.32%  '

$ws.Range("E15").Value = '  +1.35%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009584'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.18%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.101'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.99%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.129.47'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.67%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.65'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.59%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '227.96'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.43%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.00%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.195'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.07%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '160.58'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.51%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1434'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.37%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.562'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.79%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.98'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.98%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.505'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.57%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.159'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.78%  '

$ws.Range("E30").Value = '  +1.28%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05469'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.20%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.201'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.10%  '

$ws.Range("E33").Value = '  +1.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7484'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.41%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.142'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.84%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.662'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.86%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.250.54'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.24%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.762'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.12%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01792'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.16%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.675'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.94%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9054'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.94%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.001'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.07%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.00000000130'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +9.39%  '

$ws.Range("E44").Value = '  +0.06%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.983.92'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.26%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.42'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.62%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5108'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.00%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4088'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.82%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.038'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.06%  '

$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.654'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.57%  '

$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.785'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.23%  '
